$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 256 (existing rows 256-263 shift down to 257-264).
$ws.Rows.Item(256).Insert()

# Populate the newly inserted row 256 with the new weekly data point.
$ws.Cells.Item(256, 1).Value = 5
$ws.Cells.Item(256, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(256, 3).Value = "Maule"
$ws.Cells.Item(256, 4).Value = 45008
$ws.Cells.Item(256, 5).Value = 7
$ws.Cells.Item(256, 6).Value = 100112017
$ws.Cells.Item(256, 7).Value = "Apio"
$ws.Cells.Item(256, 8).Value = "Americana (o)"
$ws.Cells.Item(256, 9).Value = "Primera"
$ws.Cells.Item(256, 10).Value = 500
$ws.Cells.Item(256, 11).Value = 8000
$ws.Cells.Item(256, 12).Value = 8000
$ws.Cells.Item(256, 13).Value = 8000
$ws.Cells.Item(256, 14).Value = "`$/docena de matas"
$ws.Cells.Item(256, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(256, 16).Value = 1333
$ws.Cells.Item(256, 17).Value = 6
$ws.Cells.Item(256, 18).Value = "Hortaliza"
